# "updated sample data and SearchLine"
# Populate the sample order-items header row and pre-format the
# "Total Price" column for currency, then leave the selection where the
# user was last browsing (cell C6) and tidy column A's width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:D1) - becomes shared strings 0..3
$ws.Range("A1").Value = "OrderNo"
$ws.Range("B1").Value = "Product"
$ws.Range("C1").Value = "Fulfilment Store"
$ws.Range("D1").Value = "Total Price"

# Pre-format the (currently empty) first data cell under "Total Price"
# as currency - dollars, no decimals, red negatives.
$ws.Range("D2").NumberFormat = '"$"#,##0;[Red]\-"$"#,##0'

# Autosize the OrderNo column now that it has data in it.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Restore the last active selection.
$ws.Range("C6").Select()
